# Quarterly indexing esoteric bug-fix operation
# Column A holds quarter-end timestamps (e.g. 2010-12-31 23:59:59.999) that were
# meant to be the 15th day of the month 1.5 months earlier, at midnight
# (e.g. 2010-11-15 00:00:00). Recompute each value in column A, rows 2-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 47; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2

    $oldDate = [DateTime]::FromOADate($oldVal)
    $shifted = $oldDate.AddMonths(-1)
    $deltaDays = 15 - $shifted.Day
    $newDate = $shifted.AddDays($deltaDays).Date

    $cell.Value = $newDate.ToOADate()
}
